$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (i.e. right
#    before "总计"), matching the layout used by the other quarter sheets.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160
$newSheet.Range("B1:H1").Borders.LineStyle = 1

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "162006"
$newSheet.Range("C2").Value = "长城久富核心成长混合(LOF)"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "19.40"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "79.43"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.70"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.9118"
$newSheet.Range("H2").Value = 5

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "519967"
$newSheet.Range("C3").Value = "长信利富债券"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "5.83"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "20.15"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "1.19"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0694"
$newSheet.Range("H3").Value = 2

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "000976"
$newSheet.Range("C4").Value = "长城新兴产业灵活配置混合"
$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "1.02"
$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "77.95"
$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "4.61"
$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.0470"
$newSheet.Range("H4").Value = 5

$newSheet.Range("A2:A4").Font.Bold = $true
$newSheet.Range("A2:A4").HorizontalAlignment = -4108
$newSheet.Range("A2:A4").VerticalAlignment = -4160
$newSheet.Range("A2:A4").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row at the top of the data
#    (2022-Q1: 3 funds, 1.03 亿元) and push the existing rows down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2:D2").Insert(-4121)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.03

# Match the existing index-column style (bold, centred, thin border) by
# copying the format already used on the other index cells (A3, post-shift).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

Write-Host "edit complete"
